# Weekly data update: insert one new record at the top of the data block
# (row 137), pushing all existing records (rows 137-245) down by one row
# (to 138-246). This matches how the upstream source republishes its
# "Fruta / hortaliza, semanal" feed with the newest observation first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 137:245 down to 138:246, leaving a blank row 137 to fill in.
$ws.Rows.Item(137).Insert()

# Populate the new row 137 with the latest weekly observation.
$ws.Range("A137").Value = 8
$ws.Range("B137").Value = "Terminal La Palmera de La Serena"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44574
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 100112032
$ws.Range("G137").Value = "Zapallo italiano"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 560
$ws.Range("K137").Value = 10000
$ws.Range("L137").Value = 11000
$ws.Range("M137").Value = 10500
$ws.Range("N137").Value = "`$/caja 70 unidades"
$ws.Range("O137").Value = "Provincia de Limarí"
$ws.Range("P137").Value = 150
$ws.Range("Q137").Value = 70
$ws.Range("R137").Value = "Hortaliza"
